# Append a new case (row 12) to the "ランサーズ" sheet and push the former
# row 12 ("限定公開 限定公開の仕事") down to row 13, refreshing every
# "取得日時" timestamp in column A to 2026-01-10 01:57:12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ts = "2026-01-10 01:57:12"

# --- 1. Capture the current (pre-edit) row 12 values before we overwrite them ---
$oldB12 = $ws.Range("B12").Text
$oldC12 = $ws.Range("C12").Text
$oldD12 = $ws.Range("D12").Text
$oldE12 = $ws.Range("E12").Text
$oldF12 = $ws.Range("F12").Text
$oldG12 = $ws.Range("G12").Value2

# --- 2. Drop every hyperlink on the sheet; Excel re-links them as the cell
#        text is rewritten below (this runtime does not support editing an
#        existing Hyperlink.Address/Delete in place, so a clean rebuild is
#        the reliable path). ---
$ws.Hyperlinks.Delete()

# --- 3. Move the old row 12 down to row 13, unchanged apart from the
#        refreshed timestamp. ---
$ws.Range("A13").Value = $ts
$ws.Range("B13").Value = $oldB12
$ws.Range("C13").Value = $oldC12
$ws.Range("D13").Value = $oldD12
$ws.Range("E13").Value = $oldE12
$ws.Range("F13").Value = $oldF12
$ws.Range("G13").Value = $oldG12
$ws.Hyperlinks.Add($ws.Range("F13"), $oldF12, "", "", $oldF12)

# --- 4. Write the new case into row 12. ---
$ws.Range("A12").Value = $ts
$ws.Range("B12").Value = "【緊急対応】インターネットを活用した電話通知システム構築"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5468565"
$ws.Range("G12").Value = 33
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5468565", "", "", "https://www.lancers.jp/work/detail/5468565")

# --- 5. Re-create the untouched hyperlinks (F2:F11) that were dropped in
#        step 2, and refresh their row's timestamp. ---
$urls = @(
    "https://www.lancers.jp/work/detail/5460562",
    "https://www.lancers.jp/work/detail/5468493",
    "https://www.lancers.jp/work/detail/5468303",
    "https://www.lancers.jp/work/detail/5460563",
    "https://www.lancers.jp/work/detail/5467745",
    "https://www.lancers.jp/work/detail/5468441",
    "https://www.lancers.jp/work/detail/5467910",
    "https://www.lancers.jp/work/detail/5468432",
    "https://www.lancers.jp/work/detail/5468347",
    "https://www.lancers.jp/work/detail/5467981"
)

for ($i = 0; $i -lt $urls.Count; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $ts
    $cell = $ws.Range("F$r")
    $ws.Hyperlinks.Add($cell, $urls[$i], "", "", $urls[$i])
}
